# Mise à jour de l'application
# Add a new day's attendance column (2025-10-09) to the sheet: a new date
# header in BI1, and each player's attendance status for that day in
# column BI (same status as the previous session, except for two players
# whose status changed for this new date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header (09/10/2025 -> Excel serial 45939)
$ws.Range("BI1").Value = 45939

# Per-player attendance status for the new date column (row -> status code)
# P = Présent, A = Absent, B = Blessure, M = Malade, R = Réserve
$ws.Range("BI2").Value  = "P"
$ws.Range("BI3").Value  = "R"
$ws.Range("BI4").Value  = "P"
$ws.Range("BI5").Value  = "P"
$ws.Range("BI6").Value  = "A"
$ws.Range("BI7").Value  = "P"
$ws.Range("BI8").Value  = "B"
$ws.Range("BI9").Value  = "P"
$ws.Range("BI10").Value = "P"
$ws.Range("BI11").Value = "P"
# Row 12 (Yanis Berrached) has no entry for this date, same as before.
$ws.Range("BI13").Value = "B"
$ws.Range("BI14").Value = "P"
$ws.Range("BI15").Value = "P"
$ws.Range("BI16").Value = "P"
$ws.Range("BI17").Value = "B"
$ws.Range("BI18").Value = "P"
$ws.Range("BI19").Value = "P"
$ws.Range("BI20").Value = "M"
$ws.Range("BI21").Value = "M"
$ws.Range("BI22").Value = "P"
$ws.Range("BI23").Value = "B"
$ws.Range("BI24").Value = "P"
$ws.Range("BI25").Value = "P"
$ws.Range("BI26").Value = "P"
$ws.Range("BI27").Value = "P"
$ws.Range("BI28").Value = "P"
$ws.Range("BI29").Value = "P"

# Copy the formatting from the previous date column (BH) onto the new one
# so the new cells pick up the same style (date number format / centering)
# instead of defaulting to "General".
$ws.Range("BH1:BH11").Copy()
$ws.Range("BI1:BI11").PasteSpecial(-4122)

$ws.Range("BH13:BH29").Copy()
$ws.Range("BI13:BI29").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the active selection forward, matching where the user left off.
$ws.Range("BJ25").Select()
